# Auto-generated Excel COM-interop script to apply Marilith_Profits.xlsx data refresh
# Updates LeveProfit/Price derived columns (H-N) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 992.6667
$ws.Range("I80").Value = 830.8333
$ws.Range("J80").Value = 1316.3334
$ws.Range("K80").Value = 2492.4999
$ws.Range("L80").Value = 3949.0002
$ws.Range("M80").Value = -1494.4999
$ws.Range("N80").Value = -5945.0002

$ws.Range("H83").Value = 992.6667
$ws.Range("I83").Value = 830.8333
$ws.Range("J83").Value = 1316.3334
$ws.Range("K83").Value = 7477.4997
$ws.Range("L83").Value = 11847.0006
$ws.Range("M83").Value = -2485.4997
$ws.Range("N83").Value = -21831.0006

$ws.Range("H98").Value = 1953.091
$ws.Range("I98").Value = 1186
$ws.Range("K98").Value = 1186
$ws.Range("M98").Value = 312

$ws.Range("H113").Value = 18449.75
$ws.Range("J113").Value = 21266.334
$ws.Range("L113").Value = 21266.334
$ws.Range("N113").Value = -27774.334

$ws.Range("H122").Value = 1953.091
$ws.Range("I122").Value = 1186
$ws.Range("K122").Value = 3558
$ws.Range("M122").Value = -1108

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1985.625
$ws.Range("I2").Value = 2097.8572
$ws.Range("K2").Value = 2097.8572
$ws.Range("M2").Value = -1984.8572

$ws.Range("H32").Value = 2278.2295
$ws.Range("I32").Value = 2101.2205
$ws.Range("K32").Value = 2101.2205
$ws.Range("M32").Value = -1814.2205

$ws.Range("H61").Value = 1694
$ws.Range("J61").Value = 1398.75
$ws.Range("L61").Value = 1398.75
$ws.Range("N61").Value = -1822.75

$ws.Range("H74").Value = 1037.2
$ws.Range("I74").Value = 1029.4445
$ws.Range("K74").Value = 1029.4445
$ws.Range("M74").Value = -155.4445000000001

$ws.Range("H77").Value = 1037.2
$ws.Range("I77").Value = 1029.4445
$ws.Range("K77").Value = 5147.2225
$ws.Range("M77").Value = -779.2224999999999

$ws.Range("H97").Value = 886
$ws.Range("I97").Value = 810.5714
$ws.Range("J97").Value = 1150
$ws.Range("K97").Value = 810.5714
$ws.Range("L97").Value = 1150
$ws.Range("M97").Value = -314.5714
$ws.Range("N97").Value = -2142

$ws.Range("H116").Value = 1985.625
$ws.Range("I116").Value = 2097.8572
$ws.Range("K116").Value = 2097.8572
$ws.Range("M116").Value = 196.1428000000001

$ws.Range("H119").Value = 698
$ws.Range("J119").Value = 698
$ws.Range("L119").Value = 698
$ws.Range("N119").Value = -10374

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 3390
$ws.Range("I132").Value = 3390
$ws.Range("K132").Value = 10170
$ws.Range("M132").Value = -7640

$ws.Range("H136").Value = 1694
$ws.Range("J136").Value = 1398.75
$ws.Range("L136").Value = 4196.25
$ws.Range("N136").Value = -9296.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1985.625
$ws.Range("I3").Value = 2097.8572
$ws.Range("K3").Value = 2097.8572
$ws.Range("M3").Value = -1983.8572

$ws.Range("H86").Value = 4203.1113
$ws.Range("I86").Value = 3974
$ws.Range("J86").Value = 4563.143
$ws.Range("K86").Value = 3974
$ws.Range("L86").Value = 4563.143
$ws.Range("M86").Value = -2851
$ws.Range("N86").Value = -6809.143

$ws.Range("H89").Value = 4203.1113
$ws.Range("I89").Value = 3974
$ws.Range("J89").Value = 4563.143
$ws.Range("K89").Value = 19870
$ws.Range("L89").Value = 22815.715
$ws.Range("M89").Value = -14254
$ws.Range("N89").Value = -34047.715

$ws.Range("H105").Value = 4477.857
$ws.Range("I105").Value = 4789.4
$ws.Range("K105").Value = 4789.4
$ws.Range("M105").Value = -3042.4

$ws.Range("H107").Value = 1389.75
$ws.Range("J107").Value = 2166.6667
$ws.Range("L107").Value = 2166.6667
$ws.Range("N107").Value = -6006.6667

$ws.Range("H134").Value = 8091
$ws.Range("I134").Value = 8346.6875
$ws.Range("K134").Value = 25040.0625
$ws.Range("M134").Value = -22505.0625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3994
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 3994
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 3994
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -4400

$ws.Range("H105").Value = 519.6
$ws.Range("I105").Value = 366
$ws.Range("K105").Value = 366
$ws.Range("M105").Value = 1381

$ws.Range("H132").Value = 3399.6667
$ws.Range("I132").Value = 3399.6667
$ws.Range("K132").Value = 10199.0001
$ws.Range("M132").Value = -7669.000100000001

$ws.Range("H134").Value = 1997.25
$ws.Range("I134").Value = 1997.25
$ws.Range("K134").Value = 5991.75
$ws.Range("M134").Value = -3456.75

$ws.Range("H136").Value = 3994
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 3994
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 11982
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -17082

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 998
$ws.Range("I57").Value = 998
$ws.Range("K57").Value = 2994
$ws.Range("M57").Value = -2435

$ws.Range("H81").Value = 2666.3333
$ws.Range("J81").Value = 1500
$ws.Range("L81").Value = 4500
$ws.Range("N81").Value = -6746

$ws.Range("H84").Value = 2666.3333
$ws.Range("J84").Value = 1500
$ws.Range("L84").Value = 13500
$ws.Range("N84").Value = -24732

$ws.Range("H112").Value = 10
$ws.Range("I112").Value = 10
$ws.Range("K112").Value = 30
$ws.Range("M112").Value = 1078

$ws.Range("H115").Value = 225
$ws.Range("I115").Value = 225
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 675
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 500
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 21666.666
$ws.Range("J15").Value = 21666.666
$ws.Range("L15").Value = 21666.666
$ws.Range("N15").Value = -22242.666

$ws.Range("H57").Value = 22395.8
$ws.Range("J57").Value = 24994.75
$ws.Range("L57").Value = 24994.75
$ws.Range("N57").Value = -26634.75

$ws.Range("H81").Value = 21666.666
$ws.Range("J81").Value = 21666.666
$ws.Range("L81").Value = 21666.666
$ws.Range("N81").Value = -23662.666

$ws.Range("H84").Value = 21666.666
$ws.Range("J84").Value = 21666.666
$ws.Range("L84").Value = 64999.99800000001
$ws.Range("N84").Value = -74983.99800000001

$ws.Range("H97").Value = 792.4
$ws.Range("I97").Value = 766.125
$ws.Range("K97").Value = 766.125
$ws.Range("M97").Value = -270.125

$ws.Range("H122").Value = 4302.5454
$ws.Range("I122").Value = 4226.3335
$ws.Range("K122").Value = 12679.0005
$ws.Range("M122").Value = -10229.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1538.6
$ws.Range("I16").Value = 1538.6
$ws.Range("K16").Value = 1538.6
$ws.Range("M16").Value = -1368.6

$ws.Range("H22").Value = 1098.9231
$ws.Range("J22").Value = 1997.5
$ws.Range("L22").Value = 1997.5
$ws.Range("N22").Value = -2587.5

$ws.Range("H27").Value = 1098.9231
$ws.Range("J27").Value = 1997.5
$ws.Range("L27").Value = 1997.5
$ws.Range("N27").Value = -2211.5

$ws.Range("H40").Value = 2177.3
$ws.Range("I40").Value = 2198.111
$ws.Range("K40").Value = 2198.111
$ws.Range("M40").Value = -2062.111

$ws.Range("H46").Value = 2916.6333
$ws.Range("I46").Value = 2147
$ws.Range("K46").Value = 2147
$ws.Range("M46").Value = -1959

$ws.Range("H100").Value = 2669
$ws.Range("I100").Value = 1115
$ws.Range("K100").Value = 1115
$ws.Range("M100").Value = -574

$ws.Range("H136").Value = 3289.7778
$ws.Range("I136").Value = 3014.7144
$ws.Range("K136").Value = 9044.143199999999
$ws.Range("M136").Value = -6494.143199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3585.1904
$ws.Range("I136").Value = 3449.2896
$ws.Range("K136").Value = 10347.8688
$ws.Range("M136").Value = -7797.8688
